# Scheduled market-price refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the leve-crafting sheets with freshly scraped values.
$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$gsm = $wb.Worksheets.Item("GSM")
$ltw = $wb.Worksheets.Item("LTW")
$wvr = $wb.Worksheets.Item("WVR")

# ALC!row 86
$alc.Range("H86").Value = 4774.3477
$alc.Range("I86").Value = 1865
$alc.Range("J86").Value = 7948.1816
$alc.Range("K86").Value = 1865
$alc.Range("L86").Value = 7948.1816
$alc.Range("M86").Value = -742
$alc.Range("N86").Value = -10194.1816

# ALC!row 89
$alc.Range("H89").Value = 4774.3477
$alc.Range("I89").Value = 1865
$alc.Range("J89").Value = 7948.1816
$alc.Range("K89").Value = 9325
$alc.Range("L89").Value = 39740.908
$alc.Range("M89").Value = -3709
$alc.Range("N89").Value = -50972.908

# ALC!row 92
$alc.Range("H92").Value = 6900.625
$alc.Range("I92").Value = 13175
$alc.Range("K92").Value = 13175
$alc.Range("M92").Value = -11927

# ALC!row 107
$alc.Range("H107").Value = 921.6316
$alc.Range("I107").Value = 1151.6666
$alc.Range("K107").Value = 1151.6666
$alc.Range("M107").Value = 768.3334

# ARM!row 31
$arm.Range("H31").Value = 7656.375
$arm.Range("I31").Value = 5178.7144
$arm.Range("J31").Value = 25000
$arm.Range("K31").Value = 5178.7144
$arm.Range("L31").Value = 25000
$arm.Range("M31").Value = -4884.7144
$arm.Range("N31").Value = -25588

# ARM!row 32
$arm.Range("H32").Value = 1177113.9
$arm.Range("I32").Value = 1177113.9
$arm.Range("K32").Value = 1177113.9
$arm.Range("M32").Value = -1176826.9

# ARM!row 132
$arm.Range("H132").Value = 25023.842
$arm.Range("I132").Value = 41538.52
$arm.Range("J132").Value = 3294
$arm.Range("K132").Value = 124615.56
$arm.Range("L132").Value = 9882
$arm.Range("M132").Value = -122085.56
$arm.Range("N132").Value = -14942

# BSM!row 80
$bsm.Range("H80").Value = 454.57144
$bsm.Range("I80").Value = 597.3333
$bsm.Range("J80").Value = 347.5
$bsm.Range("K80").Value = 597.3333
$bsm.Range("L80").Value = 347.5
$bsm.Range("M80").Value = 400.6667
$bsm.Range("N80").Value = -2343.5

# BSM!row 83
$bsm.Range("H83").Value = 454.57144
$bsm.Range("I83").Value = 597.3333
$bsm.Range("J83").Value = 347.5
$bsm.Range("K83").Value = 2986.6665
$bsm.Range("L83").Value = 1737.5
$bsm.Range("M83").Value = 2005.3335
$bsm.Range("N83").Value = -11721.5

# BSM!row 102
$bsm.Range("H102").Value = 15000
$bsm.Range("I102").Value = 7000
$bsm.Range("J102").Value = 39000
$bsm.Range("K102").Value = 7000
$bsm.Range("L102").Value = 39000
$bsm.Range("M102").Value = -3755
$bsm.Range("N102").Value = -45490

# CRP!row 31
$crp.Range("H31").Value = 2117810.5
$crp.Range("I31").Value = 1112.7858
$crp.Range("J31").Value = 3139664.5
$crp.Range("K31").Value = 1112.7858
$crp.Range("L31").Value = 3139664.5
$crp.Range("M31").Value = -817.7858000000001
$crp.Range("N31").Value = -3140254.5

# CRP!row 34
$crp.Range("H34").Value = 2117810.5
$crp.Range("I34").Value = 1112.7858
$crp.Range("J34").Value = 3139664.5
$crp.Range("K34").Value = 1112.7858
$crp.Range("L34").Value = 3139664.5
$crp.Range("M34").Value = -910.7858000000001
$crp.Range("N34").Value = -3140068.5

# CRP!row 129
$crp.Range("H129").Value = 48332.668
$crp.Range("J129").Value = 48332.668
$crp.Range("L129").Value = 48332.668
$crp.Range("N129").Value = -58332.668

# CRP!row 134
$crp.Range("H134").Value = 23811062
$crp.Range("I134").Value = 38462376
$crp.Range("J134").Value = 2680.5
$crp.Range("K134").Value = 115387128
$crp.Range("L134").Value = 8041.5
$crp.Range("M134").Value = -115384593
$crp.Range("N134").Value = -13111.5

# CUL!row 68
$cul.Range("H68").Value = 887.0806
$cul.Range("I68").Value = 524.03705
$cul.Range("J68").Value = 1167.1428
$cul.Range("K68").Value = 1572.11115
$cul.Range("L68").Value = 3501.4284
$cul.Range("M68").Value = -761.1111500000002
$cul.Range("N68").Value = -5123.428400000001

# CUL!row 69
$cul.Range("H69").Value = 3231.1875
$cul.Range("I69").Value = 986
$cul.Range("J69").Value = 4251.727
$cul.Range("K69").Value = 2958
$cul.Range("L69").Value = 12755.181
$cul.Range("M69").Value = -2147
$cul.Range("N69").Value = -14377.181

# CUL!row 71
$cul.Range("H71").Value = 887.0806
$cul.Range("I71").Value = 524.03705
$cul.Range("J71").Value = 1167.1428
$cul.Range("K71").Value = 4716.33345
$cul.Range("L71").Value = 10504.2852
$cul.Range("M71").Value = -660.3334500000001
$cul.Range("N71").Value = -18616.2852

# CUL!row 72
$cul.Range("H72").Value = 3231.1875
$cul.Range("I72").Value = 986
$cul.Range("J72").Value = 4251.727
$cul.Range("K72").Value = 8874
$cul.Range("L72").Value = 38265.543
$cul.Range("M72").Value = -4818
$cul.Range("N72").Value = -46377.543

# CUL!row 107
$cul.Range("H107").Value = 1088.2903
$cul.Range("I107").Value = 490.18182
$cul.Range("J107").Value = 1417.25
$cul.Range("K107").Value = 1470.54546
$cul.Range("L107").Value = 4251.75
$cul.Range("M107").Value = 449.45454
$cul.Range("N107").Value = -8091.75

# GSM!row 15
$gsm.Range("H15").Value = 5500
$gsm.Range("J15").Value = 5500
$gsm.Range("L15").Value = 5500
$gsm.Range("N15").Value = -6076

# GSM!row 81
$gsm.Range("H81").Value = 5500
$gsm.Range("J81").Value = 5500
$gsm.Range("L81").Value = 5500
$gsm.Range("N81").Value = -7496

# GSM!row 84
$gsm.Range("H84").Value = 5500
$gsm.Range("J84").Value = 5500
$gsm.Range("L84").Value = 16500
$gsm.Range("N84").Value = -26484

# GSM!row 132
$gsm.Range("H132").Value = 2407107.5
$gsm.Range("I132").Value = 4313450
$gsm.Range("J132").Value = 3458.087
$gsm.Range("K132").Value = 12940350
$gsm.Range("L132").Value = 10374.261
$gsm.Range("M132").Value = -12937820
$gsm.Range("N132").Value = -15434.261

# LTW!row 16
$ltw.Range("H16").Value = 1888.7142
$ltw.Range("I16").Value = 1888.7142
$ltw.Range("J16").Value = 0
$ltw.Range("K16").Value = 1888.7142
$ltw.Range("L16").Value = 0
$ltw.Range("M16").Value = -1718.7142
$ltw.Range("N16").ClearContents() | Out-Null

# WVR!row 62
$wvr.Range("H62").Value = 4408.3335
$wvr.Range("I62").Value = 0
$wvr.Range("J62").Value = 4408.3335
$wvr.Range("K62").Value = 0
$wvr.Range("L62").Value = 4408.3335
$wvr.Range("M62").ClearContents() | Out-Null
$wvr.Range("N62").Value = -5656.3335

# WVR!row 65
$wvr.Range("H65").Value = 4408.3335
$wvr.Range("I65").Value = 0
$wvr.Range("J65").Value = 4408.3335
$wvr.Range("K65").Value = 0
$wvr.Range("L65").Value = 22041.6675
$wvr.Range("M65").ClearContents() | Out-Null
$wvr.Range("N65").Value = -28281.6675
